$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("CONTAS FAPERS")
$ws2 = $wb.Worksheets.Item("HP FAPERS")

# Update the HP FAPERS sheet's D-column descriptions ("autocomplete" consolidation of
# FEV/21-dated entries onto the shared, undated wording used going forward).

$ws2.Range("D2").Value = 'VLR GRATIF RICARDO ALTAIR SCHWARZ '
$ws2.Range("D3").Value = 'VLR GRATIFJOSE PEDRO OSORIO MENDINA '
$ws2.Range("D4").Value = 'VLR GRATIF ALVARO ROQUE KERN JUNQUEIRA '
$ws2.Range("D5").Value = 'VLR AJUSTE REF GRATIF E REPRESENTAÇÃO N/MES'
$ws2.Range("D6").Value = 'LANÇAMENTO PARA ZERAR O CENTRO DE CUSTEIO COMUM: 2/28/2021'
$ws2.Range("D7").Value = 'TRANSFERÊNCIA DOS SALDOS PARA PLANOS - 28/02/2021'
$ws2.Range("D8").Value = 'APROPR DE INSS RICARDO ALTAIR SCHWARZ '
$ws2.Range("D9").Value = 'APROPR DE INSS JOSE PEDRO OSORIO MENDINA '
$ws2.Range("D10").Value = 'APROPR DE INSS ALVARO ROQUE KERN JUNQUEIRA '
$ws2.Range("D11").Value = 'APROPR SEGURO DIRIGENTES  CFE APOLICE AIG SEGUROS BRASIL S/A Nº 08737.2020.01.0310.001161'
$ws2.Range("D12").Value = 'APROPR SALÁRIOS CFE FOLHA PGTO FEV/21'
$ws2.Range("D13").Value = 'VLR AJUSTE S/IRRF TRAB ASSALARIADO N/MÊS'
$ws2.Range("D14").Value = 'APROP INSS CFE FOLHA SALARIAL '
$ws2.Range("D15").Value = 'APROP FGTS CFE FOLHA SALARIAL '
$ws2.Range("D16").Value = 'APROP FGTS CFE FOLHA SALARIAL '
$ws2.Range("D17").Value = 'APROPR CONTR ASSISITIDOS PLANO MISTO PATROC FAPERS N/MES'
$ws2.Range("D18").Value = 'APROPR DE CONTR ASSISITIDOS PGS PATROC FAPERS N/MES'
$ws2.Range("D19").Value = 'APROPR DE CONTR ATIVOS PGS PATROC FAPERS N/MES'
$ws2.Range("D20").Value = 'APROPR DE CONTRIB ATIVOS PREVER PATROC FAPERS N/MES'
$ws2.Range("D21").Value = 'APROPR CONTRIB ATIVOS PLANO MISTO PATROC FAPERS N/MES'
$ws2.Range("D22").Value = 'ADM-01705/2020 -  EMIT ORC: ADM-01706/2020 - 02/02/2021-2100024502-17 - PREST SERVICOS DE CONSULTAS MEDICAS 35 - PLANO DE SAUDE - PATRONAL - UNIMED PORTO ALEGRE SOCIEDADE COOPERATIVA DE TRABALHO MEDICODT COMPET: 02/02/2021'
$ws2.Range("D23").Value = 'ADM-01714/2020 -  EMIT ORC: ADM-01715/2020 - 08/02/2021-91359-17 - PREST SERV ODONTOLOGICOS 34 - PLANO ODONTOLOGICO - PATRONAL - UNIODONTO PORTO ALEGRE - COOPERATIVA ODONTOLOGICA LTDADT COMPET: 08/02/2021'
$ws2.Range("D24").Value = 'ADM-00015/2021 -  EMIT ORC: ADM-00016/2021 - 22/02/2021-100043171679-29 - SEGURO EMPRESARIAL PARA EMPREGADOS 51 - SEGURO FUNCIONARIOS - SOMPO SEGUROS S/ADT COMPET: 22/02/2021'
$ws2.Range("D25").Value = 'APROP AUXILIO EDUCAÇÃO CFE FL SALARIAL FEV/21.'
$ws2.Range("D26").Value = 'VLR DESC VR CFE FL SALARIAL FEV/21'
$ws2.Range("D27").Value = 'APROPR SEGURO EMPREGADOS  CFE APOLICE AIG SEGUROS BRASIL S/A Nº 08737.2020.01.0310.001161'
$ws2.Range("D28").Value = 'APROPR TICKET SERV S/A REF ALIMENTAÇÃO '
$ws2.Range("D29").Value = 'APROPR TICKET SERV S/A REF REFEIÇÃO '
$ws2.Range("D30").Value = 'APROPR TICKET SERV S/A REF REFEIÇÃO ESTAGIÁRIOS '
$ws2.Range("D31").Value = 'APROPR PROV FERIAS N/MES'
$ws2.Range("D32").Value = 'APROPR PROV INSS S/FERIAS N/MES'
$ws2.Range("D33").Value = 'APROPR PROV FGTS S/FERIAS N/MES'
$ws2.Range("D34").Value = 'APROPR PROV CONTR FAPERS S/FERIAS N/MES'
$ws2.Range("D35").Value = 'APROPR PROV 13º SALARIO N/MES'
$ws2.Range("D36").Value = 'APROPR PROV INSS S/13º SALARIO N/MES'
$ws2.Range("D37").Value = 'APROPR PROV CONTR FAPERS S/13º SALÁRIO N/MES'
$ws2.Range("D38").Value = 'APROPR PROV FGTS S/13º SALÁRIO N/MES'
$ws2.Range("D39").Value = 'VLR AJUSTE REF PROVISÃO DE FÉRIAS N/MÊS'
$ws2.Range("D40").Value = 'VLR AJUSTE REF PROVISÃO 13º SALÁRIO N/MÊS'
$ws2.Range("D41").Value = 'APROPR BOLSA ESTAGIO CFE FOLHA PGTO FEV/21'
$ws2.Range("D42").Value = 'APROPR CFE REC 01342766 ANS GRAFICA REF SERVIÇOS DE IMPRESSÕES DE CARTÃO DE VISITA'
$ws2.Range("D43").Value = 'APROP AUXILIO EDUCAÇÃO CFE FL SALARIAL FEV/21'
$ws2.Range("D45").Value = 'VLR DESC CCA PARTIC PGS CFE FL PGTO FAPERS '
$ws2.Range("D46").Value = 'VLR DESC UNIMED CFE FL SALARIAL FEV/21'
$ws2.Range("D47").Value = 'VLR DESC UNIODONTO CFE FL SALARIAL FEV/21'
$ws2.Range("D48").Value = 'VLR PENSAO ALIMENTICIA CFE FL PAGTO '
$ws2.Range("D49").Value = 'VLR RET CONTR FAPERS CFE FL SALARIAL FEV/21'
$ws2.Range("D50").Value = 'VLR RET INSS CFE FL PGTO FEV/21'
$ws2.Range("D51").Value = 'VLR RET IRRF CFE FL PGTO FEV/21'
$ws2.Range("D52").Value = 'VLR TAFIC PBD I '
$ws2.Range("D53").Value = 'VLR TAFIC PGS '
$ws2.Range("D54").Value = 'VLR TAFIC PLANO MISTO '
$ws2.Range("D55").Value = 'VLR TAFIC PREVER '
$ws2.Range("D56").Value = 'APROPR VALE TRANSPORTE ESTAGIARIOS '
$ws2.Range("D57").Value = 'APROPR ALUGUEL SEDE FAPERS CASA 1073 NA RUA MARCILIO DIAS N/ MES'
$ws2.Range("D58").Value = 'VLR BX PROV CONTRIB FAPERS S/FÉRIAS N/MÊS'
$ws2.Range("D59").Value = 'VLR BX PROV FERIAS DE LIGIAMAR FRIZZO CFE REC'
$ws2.Range("D60").Value = 'VLR BX PROV FGTS S/FÉRIAS N/MÊS'
$ws2.Range("D61").Value = 'VLR BX PROV INSS S/FÉRIAS N/MÊS'
$ws2.Range("D62").Value = 'VLR CONTR FAPERS S/ FÉRIAS DE LIGIAMAR FRIZZO CFE REC'
$ws2.Range("D63").Value = 'VLR RET INSS S/ FÉRIAS DE LIGIAMAR FRIZZO CFE REC'
$ws2.Range("D64").Value = 'APROPR INSS S/ GRATIF RICARDO ALTAIR SCHWARZ '
$ws2.Range("D65").Value = 'VLR IRRF S/ GRATIF ALVARO ROQUE KERN JUNQUEIRA '
$ws2.Range("D66").Value = 'VLR IRRF S/ GRATIF JOSE PEDRO OSORIO MENDINA '
$ws2.Range("D67").Value = 'VLR IRRF S/ GRATIF RICARDO ALTAIR SCHWARZ '
$ws2.Range("D68").Value = 'VLR PARC 02/12 REF IPTU 2021 CASA 1073 MARCILIO DIAS'
$ws2.Range("D69").Value = 'APROPR SEGURO CONSELHO DELIBERATIVO  CFE APOLICE AIG SEGUROS BRASIL S/A Nº 08737.2020.01.0310.001161'
$ws2.Range("D70").Value = 'APROPR SEGURO CONSELHO FISCAL  CFE APOLICE AIG SEGUROS BRASIL S/A Nº 08737.2020.01.0310.001161'
$ws2.Range("D71").Value = 'DESC PARCELA EMPRÉST CFE FL SALARIAL FAPERS N/MÊS'
$ws2.Range("D72").Value = 'ADM-00001/2021 - 19/02/2021-2021/91-19 - ELABORACAO DA FOLHA DE PAGAMENTO E OBRIGACOES SOCIAIS 60, CFE NFSE N  2021-91 - FOLHA PAGAMENTO E CONSULT EM RH - MORESCO CONTABILIDADE S/SDT COMPET: 19/02/2021'
$ws2.Range("D73").Value = 'ADM-00004/2021 - 24/02/2021 - DESP DESLOCAMENTO AO BANRISUL RETIRADA CARTAO ASSINATURAS DIRETORIA - CONDUCAO, TRANSPORTE E TELEENTREGA - LUANA RODRIGUES DE OLIVEIRADT COMPET: 24/02/2021'
$ws2.Range("D74").Value = 'ADM-00005/2021 - 23/02/2021-3167/3168-19 - SERVICOS DE ARMAZENAGEM/DIGITALIZACAO DE DOCUMENTOS 57, CFE NFSE-N  3167 - GUARDA CONSERVACAO DE DOCUMENTOS - METROFILE BRASIL GESTAO DA INFORMACAO LTDADT COMPET: 26/02/2021'
$ws2.Range("D75").Value = 'ADM-00005/2021 - 23/02/2021-3167/3168-19 - SERVICOS DE ARMAZENAGEM/DIGITALIZACAO DE DOCUMENTOS 57, CFE NFSE-N 3168 - DIGITALIZACAO DE DOCUMENTOS - METROFILE BRASIL GESTAO DA INFORMACAO LTDADT COMPET: 26/02/2021'
$ws2.Range("D76").Value = 'ADM-00015/2021 - 22/02/2021-100043171679-29 - SEGURO EMPRESARIAL PARA EMPREGADOS 51 - SEGURO FUNCIONARIOS - SOMPO SEGUROS S/ADT COMPET: 22/02/2021'
$ws2.Range("D77").Value = 'ADM-00017/2021 - 24/02/2021-448430 - PRESTAÇÃO DE CONTAS-6 - DESLOCAMENTO DE TAXI RESIDENCIA ATE SEDE FAPERS-RECIBO EPTC N.448430 - CONDUCAO, TRANSPORTE E TELEENTREGA - ALEXANDRE VILLELA ILHADT COMPET: 22/02/2021'
$ws2.Range("D78").Value = 'ADM-00018/2021 - 23/02/2021-197/31958-6 - MAGAZINE LUIZA NFC-E 516610 E 031958 - CARGA CHIC NUCLEO ADM FINANC PARA EMPRESTIMOS - SERVICOS DE TELEFONIA - DT COMPET: 23/02/2021'
$ws2.Range("D79").Value = 'ADM-00019/2021 - CLARO RS - CARGA CHIC GESTAO DOS PLANOS - SERVICOS DE TELEFONIA - DT COMPET: 22/02/2021'
$ws2.Range("D80").Value = 'ADM-00020/2021 - 01/03/2021-BRA21X03337-6 - SERV TRANSPORTE PASSAG - CONTRATO CABIFY 26 - SERVICOS DE TRANSPORTE PASSAGEIROS - CABIFY AGENCIA DE SERVICOS DE TRANSPORTE DE PASSAGEIROS LTDADT COMPET: 26/02/2021'
$ws2.Range("D81").Value = 'ADM-00022/2021 - 25/02/2021-21/006104595-17 - SERV. AGUA E ESGOTO-CONTRATO DMAE 17 - DMAE - DEPARTAMENTO MUNICIPAL DE AGUA E ESGOTOSDT COMPET: 25/02/2021'
$ws2.Range("D82").Value = 'ADM-00024/2021 - 27/02/2021-8207-17 - SERV. POSTAIS-CONTRATO FILIZOLA 18 - DESPESAS POSTAGEM E CARTORIO - FILIZOLA CENTRO DE SERVICOS LTDADT COMPET: 27/02/2021'
$ws2.Range("D83").Value = 'ADM-00026/2021 - 25/02/2021-2021/264-19 - SERV. GRAFICOS - CONTRATO POSTAL MKT 19 - SERVICOS GRAFICOS/IMPRESSAO - POSTAL MKT INFORMATICA LTDADT COMPET: 25/02/2021'
$ws2.Range("D84").Value = 'ADM-00029/2021 - 26/02/2021-139/2021-6 - SERV PREST- REEMB DESPESAS POSTAGEM E CARTORIO 2 - DESPESAS POSTAGEM E CARTORIO - BARCELLOS ADVOCACIA EMPRESARIALDT COMPET: 26/02/2021'
$ws2.Range("D85").Value = 'ADM-00031/2021 - 19/02/2021-1230320448-0-17 - SERV VIVO FIXO - CONTRATO VIVO TELEFONICA 26 - SERVICOS DE TELEFONIA - TELEFONICA BRASIL S/A - VIVODT COMPET: 19/02/2021'
$ws2.Range("D86").Value = 'ADM-00033/2021 - 26/02/2021-9311-19 - SERV. LOCACAO IMPRESSORAS/MULTIFUNC-CONTRATO PRINTMAX 12 - SERVICOS DE INFORMATICA/INTERNET - PRINTMAX GERENCIAMENTO DE DOCUMENTOS EIRELIDT COMPET: 26/02/2021'
$ws2.Range("D87").Value = 'ADM-00035/2021 - 17/02/2021-36077-19 - SERVICOS DE LIMPEZA PREDIAL-CONTRATO SULCLEAN 5 - INSS S/ SERVICO - SULCLEAN SERVICOS LTDADT COMPET: 17/02/2021'
$ws2.Range("D88").Value = 'ADM-00035/2021 - 17/02/2021-36077-19 - SERVICOS DE LIMPEZA PREDIAL-CONTRATO SULCLEAN 5 - IRRF S/ SERV PRESTADO - SULCLEAN SERVICOS LTDADT COMPET: 17/02/2021'
$ws2.Range("D89").Value = 'ADM-00035/2021 - 17/02/2021-36077-19 - SERVICOS DE LIMPEZA PREDIAL-CONTRATO SULCLEAN 5 - ISSQN - SULCLEAN SERVICOS LTDADT COMPET: 17/02/2021'
$ws2.Range("D90").Value = 'ADM-00035/2021 - 17/02/2021-36077-19 - SERVICOS DE LIMPEZA PREDIAL-CONTRATO SULCLEAN 5 - SERVICOS DE LIMPEZA - SULCLEAN SERVICOS LTDADT COMPET: 17/02/2021'
$ws2.Range("D91").Value = 'ADM-00037/2021 - 01/03/2021-2021/47-19 - ENTREGA DE MATERIAIS ATRAVES DE MOTOBOY - CONDUCAO, TRANSPORTE E TELEENTREGA - MAXIMO MOTOBOY EXPRESSO LTDADT COMPET: 26/02/2021'
$ws2.Range("D92").Value = 'ADM-00040/2021 - 02/03/2021-2021/134-19 - CONSULT, MANUT E SERV TECNICOS EM INFORMTICA 15 - SERVICOS DE INFORMATICA/INTERNET - WF TECNOLOGIA LTDADT COMPET: 26/02/2021'
$ws2.Range("D93").Value = 'ADM-00047/2021 - 03/03/2021-0385374102-17 - SERV MOVEL CELULARES - CONTRATO VIVO TELEFONICA 25 - SERVICOS DE TELEFONIA - TELEFONICA BRASIL S/A - VIVODT COMPET: 26/02/2021'
$ws2.Range("D94").Value = 'ADM-00054/2021 - 09/03/2021-NFS- E 2021/30-19 - CONSERTO DE FECHADURA E COPIA DE  CHAVE TETRA -  VINICIUS RAMOS KRETZKEI'
$ws2.Range("D95").Value = 'ADM-00056/2021 -  EMIT ORC: ADM-00057/2021 - 19/02/2021-27278534-19 - CONCESSAO DO PAT-VALE RANCHO 40 - VALE ALIMENTACAO - PATRONAL - TICKET SERVICOS S.A.DT COMPET: 19/02/2021'
$ws2.Range("D96").Value = 'ADM-00056/2021 - 19/02/2021-27278534-19 - CONCESSAO DO PAT-VALE RANCHO 40 - VALE ALIMENTACAO - PATRONAL - TICKET SERVICOS S.A.DT COMPET: 19/02/2021'
$ws2.Range("D97").Value = 'ADM-00058/2021 -  EMIT ORC: ADM-00059/2021 - 19/02/2021-27278485-19 - CONCESSAO PAT- VALE REFEICAO 41 - VALE REFEICAO-PATRONAL - TICKET SERVICOS S.A.DT COMPET: 19/02/2021'
$ws2.Range("D98").Value = 'ADM-00058/2021 - 19/02/2021-27278485-19 - CONCESSAO PAT- VALE REFEICAO 41 - VALE REFEICAO-PATRONAL - TICKET SERVICOS S.A.DT COMPET: 19/02/2021'
$ws2.Range("D100").Value = 'ADM-00060/2021 - 26/02/2021-10952598-17 - CONCESSAO TRI - VALE TRASPORTE 42 - VALE TRANSPORTE - ASSOC EMPRESAS TRANSPORTES PASSAGEIROS DE POADT COMPET: 26/02/2021'
$ws2.Range("D101").Value = 'ADM-00062/2021 - 25/02/2021-1929570-6 -  4º  TABELIONATO DE NOTAS -  REC.1729570 - DESPESAS DE AUTENTICAÇÕES - DT COMPET: 24/02/2021'
$ws2.Range("D102").Value = 'ADM-00063/2021 - 01/03/2021-23128330-17 - CONVENIO AGENTE DE INTEGRACAO BOLSA DE ESTAGIO 54 - PGTO TX ADM BOLSA ESTAGIO - CIEERS - CENTRO DE INTEGRACAO EMPRESA ESCOLADT COMPET: 26/02/2021'
$ws2.Range("D103").Value = 'ADM-01621/2020 - 01/02/2021 - DESP.AUTENTICACOES 4.TABELIONATO REC. N.1726460 - DESPESAS POSTAGEM E CARTORIO - LUANA RODRIGUES DE OLIVEIRADT COMPET: 01/02/2021'
$ws2.Range("D104").Value = 'ADM-01621/2020 - 01/02/2021 - DESPESAS DESLOCAMENTO UBER - CONDUCAO, TRANSPORTE E TELEENTREGA - LUANA RODRIGUES DE OLIVEIRADT COMPET: 01/02/2021'
$ws2.Range("D105").Value = 'ADM-01646/2020 - 01/02/2021-00022031-19 - LICENCA SOFTWARE QUANTUM INDICADORES INVESTIMENTOS 37 - CONSULTORIA INVESTIMENTOS - PARAMITA TECNOLOGIA CONSULTORIA FINANCEIRA LTDADT COMPET: 01/02/2021'
$ws2.Range("D106").Value = 'ADM-01646/2020 - 01/02/2021-00022031-19 - LICENCA SOFTWARE QUANTUM INDICADORES INVESTIMENTOS 37 - IRRF S/ SERV PRESTADO - PARAMITA TECNOLOGIA CONSULTORIA FINANCEIRA LTDADT COMPET: 01/02/2021'
$ws2.Range("D107").Value = 'ADM-01646/2020 - 01/02/2021-00022031-19 - LICENCA SOFTWARE QUANTUM INDICADORES INVESTIMENTOS 37 - PIS/COFINS/CSLL SERV PRESTADO - PARAMITA TECNOLOGIA CONSULTORIA FINANCEIRA LTDADT COMPET: 01/02/2021'
$ws2.Range("D108").Value = 'ADM-01650/2020 - 10/02/2021-V80018254902-29 - CONTRIB ASSOCIATIVA - CONTRATO ABRAPP 20 - ENTIDADES E ASSOCIACOES - ASSOCIACAO BRAS. DE ENTIDADES FECHADAS PREVIDENCIA PRIVADADT COMPET: 01/02/2021'
$ws2.Range("D109").Value = 'ADM-01680/2020 - 04/02/2021-2017/17-19 - ASSES. E CONSUL. DE INVESTIMENTOS-CONTRATO MARCO & MARCO 6 - CONSULTORIA INVESTIMENTOS - MARCO & MARCO CONSULTORES FINANCEIROS ASSOCIADOS LTDADT COMPET: 04/02/2021'
$ws2.Range("D110").Value = 'ADM-01680/2020 - 04/02/2021-2017/17-19 - ASSES. E CONSUL. DE INVESTIMENTOS-CONTRATO MARCO & MARCO 6 - IRRF S/ SERV PRESTADO - MARCO & MARCO CONSULTORES FINANCEIROS ASSOCIADOS LTDADT COMPET: 04/02/2021'
$ws2.Range("D111").Value = 'ADM-01680/2020 - 04/02/2021-2017/17-19 - ASSES. E CONSUL. DE INVESTIMENTOS-CONTRATO MARCO & MARCO 6 - PIS/COFINS/CSLL SERV PRESTADO - MARCO & MARCO CONSULTORES FINANCEIROS ASSOCIADOS LTDADT COMPET: 04/02/2021'
$ws2.Range("D112").Value = 'ADM-01687/2020 - 05/02/2021-80510067-29 - SERV. HOSPEDAGEM DE SITES-CONTRATO KINGHOST 14 - SERVICOS DE INFORMATICA/INTERNET - CYBERWEB NETWORKS LTDADT COMPET: 05/02/2021'
$ws2.Range("D113").Value = 'ADM-01699/2020 - 08/02/2021-PRESTACAO CONTA-6 - TRANSPORTE DE UBER RESIDENCIA ATE FAPERS - CONDUCAO, TRANSPORTE E TELEENTREGA - CASSIO ZARPELONDT COMPET: 08/02/2021'
$ws2.Range("D114").Value = 'ADM-01701/2020 - 08/02/2021-PRESTACAO CONTA-6 - PGTO EMOLUMENTOS NO CARTORIO DE REGISTRO ESPECIAL- DIFERENCA - DESPESAS POSTAGEM E CARTORIO - CATIA BERGER ROLIMDT COMPET: 08/02/2021'
$ws2.Range("D115").Value = 'ADM-01702/2020 - 09/02/2021-PRESTACAO CONTA-6 - DELOCAMENTO DE ATA VIA RODOVIARIA - CONDUCAO, TRANSPORTE E TELEENTREGA - DIRLEI MATOS DE SOUZADT COMPET: 03/02/2021'
$ws2.Range("D116").Value = 'ADM-01703/2020 -  EMIT ORC: ADM-01704/2020 - 02/02/2021-2100024502-29 - PREST SERVICOS MEDICOS-CO-PARTICIPACAO 36 - PLANO DE SAUDE - EMPREGADOS - UNIMED PORTO ALEGRE SOCIEDADE COOPERATIVA DE TRABALHO MEDICODT COMPET: 02/02/2021'
$ws2.Range("D117").Value = 'ADM-01703/2020 - 02/02/2021-2100024502-29 - PREST SERVICOS MEDICOS-CO-PARTICIPACAO 36 - PLANO DE SAUDE - EMPREGADOS - UNIMED PORTO ALEGRE SOCIEDADE COOPERATIVA DE TRABALHO MEDICODT COMPET: 02/02/2021'
$ws2.Range("D118").Value = 'ADM-01705/2020 -  EMIT ORC: ADM-01707/2020 - 02/02/2021-2100024502-17 - PREST SERVICOS DE CONSULTAS MEDICAS 35 - PLANO DE SAUDE - EMPREGADOS - UNIMED PORTO ALEGRE SOCIEDADE COOPERATIVA DE TRABALHO MEDICODT COMPET: 02/02/2021'
$ws2.Range("D119").Value = 'ADM-01705/2020 - 02/02/2021-2100024502-17 - PREST SERVICOS DE CONSULTAS MEDICAS 35 - PLANO DE SAUDE - EMPREGADOS - UNIMED PORTO ALEGRE SOCIEDADE COOPERATIVA DE TRABALHO MEDICODT COMPET: 02/02/2021'
$ws2.Range("D120").Value = 'ADM-01705/2020 - 02/02/2021-2100024502-17 - PREST SERVICOS DE CONSULTAS MEDICAS 35 - PLANO DE SAUDE - PATRONAL - UNIMED PORTO ALEGRE SOCIEDADE COOPERATIVA DE TRABALHO MEDICODT COMPET: 02/02/2021'
$ws2.Range("D121").Value = 'ADM-01714/2020 -  EMIT ORC: ADM-01716/2020 - 08/02/2021-91359-17 - PREST SERV ODONTOLOGICOS 34 - PLANO ODONTOLOGICO - EMPREGADOS - UNIODONTO PORTO ALEGRE - COOPERATIVA ODONTOLOGICA LTDADT COMPET: 08/02/2021'
$ws2.Range("D122").Value = 'ADM-01714/2020 - 08/02/2021-91359-17 - PREST SERV ODONTOLOGICOS 34 - PLANO ODONTOLOGICO - EMPREGADOS - UNIODONTO PORTO ALEGRE - COOPERATIVA ODONTOLOGICA LTDADT COMPET: 08/02/2021'
$ws2.Range("D123").Value = 'ADM-01714/2020 - 08/02/2021-91359-17 - PREST SERV ODONTOLOGICOS 34 - PLANO ODONTOLOGICO - PATRONAL - UNIODONTO PORTO ALEGRE - COOPERATIVA ODONTOLOGICA LTDADT COMPET: 08/02/2021'
$ws2.Range("D124").Value = 'ADM-01717/2020 - 08/02/2021-PRESTACAO CONTA-6 - PAGAMENTO DE DESPESAS DE CARTORIO - DESPESAS POSTAGEM E CARTORIO - CATIA BERGER ROLIMDT COMPET: 08/02/2021'
$ws2.Range("D125").Value = 'ADM-01717/2020 - 08/02/2021-PRESTACAO CONTA-6 - PGTO TAXIS E APLICATIVOS PARA DESCOLAMENTO AO CARTORIO - CONDUCAO, TRANSPORTE E TELEENTREGA - CATIA BERGER ROLIMDT COMPET: 08/02/2021'
$ws2.Range("D126").Value = 'ADM-01718/2020 - 08/02/2021 - 4º TABELIONATO DE NOTAS - DESPESAS DE XEROX E AUTENTICACOES CFE.RECIBO N  1727468 - DESPESAS POSTAGEM E CARTORIO - DT COMPET: 08/02/2021'
$ws2.Range("D127").Value = 'ADM-01720/2020 - 04/02/2021-0385374102-29 - SERV MOVEL CELULARES - CONTRATO VIVO TELEFONICA 25 - SERVICOS DE TELEFONIA - TELEFONICA BRASIL S/A - VIVODT COMPET: 04/02/2021'
$ws2.Range("D128").Value = 'ADM-01726/2020 - 01/03/2021-3069-19 - LICENCA DE USO DE SOFTWARE-CONTRATO INTECH 1 - ASSES/CONSULT INFORMATICA - INTECH SOLUCOES EM TECNOLOGIA DA INFORMACAO LTDADT COMPET: 26/02/2021'
$ws2.Range("D129").Value = 'ADM-01726/2020 - 01/03/2021-3069-19 - LICENCA DE USO DE SOFTWARE-CONTRATO INTECH 1 - IRRF S/ SERV PRESTADO - INTECH SOLUCOES EM TECNOLOGIA DA INFORMACAO LTDADT COMPET: 26/02/2021'
$ws2.Range("D130").Value = 'ADM-01728/2020 - 26/02/2021-2021/51-19 - SERV. ADVOC. E DE CONSULTORIA JURIDICA-CONTRATO BARCELLOS 2 - CONSULTORIA JURIDICA - BARCELLOS ADVOCACIA EMPRESARIALDT COMPET: 26/02/2021'
$ws2.Range("D131").Value = 'ADM-01728/2020 - 26/02/2021-2021/51-19 - SERV. ADVOC. E DE CONSULTORIA JURIDICA-CONTRATO BARCELLOS 2 - IRRF S/ SERV PRESTADO - BARCELLOS ADVOCACIA EMPRESARIALDT COMPET: 26/02/2021'
$ws2.Range("D132").Value = 'ADM-01732/2020 - 26/02/2021-15993021-29 - SERV. E-MAIL MARKTING-CONTRATO KINGHOST 14 - SERVICOS DE INFORMATICA/INTERNET - CYBERWEB NETWORKS LTDADT COMPET: 26/02/2021'
$ws2.Range("D133").Value = 'ADM-01736/2020 - 20/02/2021-3575770-29 - SERV NET - CONTRATO VIVO TELEFONICA 22 - SERVICOS DE TELEFONIA - TELEFONICA BRASIL S/A - VIVODT COMPET: 26/02/2021'
$ws2.Range("D134").Value = 'ADM-01738/2020 - 26/02/2021-999979493598-29 - SERV VIVO 0800 - CONTRATO VIVO TELEFONICA 23 - SERVICOS DE TELEFONIA - TELEFONICA BRASIL S/A - VIVODT COMPET: 26/02/2021'
$ws2.Range("D135").Value = 'ADM-01740/2020 - 01/03/2021-1238934203-0-17 - SERV VIVO INTERNET- CONTRATO VIVO TELEFONICA 24 - SERVICOS DE TELEFONIA - TELEFONICA BRASIL S/A - VIVODT COMPET: 26/02/2021'
$ws2.Range("D136").Value = 'ADM-01744/2020 - 02/03/2021-2021/26-19 - PREST SERV MAPEAMENTO,OTIMIZ PROCESSOS 50 - ASSES/CONSULT MAPEAM,OTIMIZ PROCESSOS - MENTHOR CONSULTORIA E CAPACITACAO EMPRESARIAL LTDADT COMPET: 26/02/2021'
$ws2.Range("D137").Value = 'ADM-01746/2020 - 03/03/2021-2021/13-19 - CONSULTORIA ATUARIAL-METODO ATUARIAL 58 - CONSULTORIA ATUARIAL - METODO ATUARIAL SOCIEDADE SIMPLES LTDADT COMPET: 26/02/2021'
$ws2.Range("D138").Value = 'ADM-01746/2020 - 03/03/2021-2021/13-19 - CONSULTORIA ATUARIAL-METODO ATUARIAL 58 - IRRF S/ SERV PRESTADO - METODO ATUARIAL SOCIEDADE SIMPLES LTDADT COMPET: 26/02/2021'
$ws2.Range("D139").Value = 'ADM-01748/2020 - 01/03/2021-43-19 - CONSULTORIA FINANCEIRA E DE INVESTIMENTOS 59 - CONSULTORIA INVESTIMENTOS - C R NEUENFELDT MEDT COMPET: 26/02/2021'
$ws2.Range("D140").Value = 'ADM-01750/2020 - 04/03/2021-2021/124-19 - ASSES/CONSULT-SUPORTE TECNICO SITE FAPERS 34 - IRRF S/ SERV PRESTADO - MIRADOR ASSESSORIA ATUARIAL LTDADT COMPET: 26/02/2021'
$ws2.Range("D141").Value = 'ADM-01750/2020 - 04/03/2021-2021/124-19 - ASSES/CONSULT-SUPORTE TECNICO SITE FAPERS 34 - SERVICOS TECNICO SUPORTE SITE - MIRADOR ASSESSORIA ATUARIAL LTDADT COMPET: 26/02/2021'
$ws2.Range("D142").Value = 'ADM-01752/2020 - 04/03/2021-2021/125-19 - ASSESS/CONSULT - SERV AREA COMUNICACAO/REDE SOCIAIS 33 - ASSES/CONSULT GESTAO/PLANEJ. ESTRATEGICO - MIRADOR ASSESSORIA ATUARIAL LTDADT COMPET: 26/02/2021'
$ws2.Range("D143").Value = 'ADM-01752/2020 - 04/03/2021-2021/125-19 - ASSESS/CONSULT - SERV AREA COMUNICACAO/REDE SOCIAIS 33 - IRRF S/ SERV PRESTADO - MIRADOR ASSESSORIA ATUARIAL LTDADT COMPET: 26/02/2021'
$ws2.Range("D144").Value = 'ADM-01754/2020 - 10/02/2021-016970-19 - SERV. TECNICO SENSOR ALARME-CONTRATO RUDDER 11 - SERVICOS E MAT DE SEGURANCA/VIGILANCIA - RUDDER EQUIPAMENTOS SISTEMAS SEGURANCA LTDADT COMPET: 10/02/2021'
$ws2.Range("D145").Value = 'ADM-01756/2020 - 10/02/2021-2021/2206-19 - SERV. MONITORAMENTO-TELEALARME-CONTRATO RUDDER 11 - SERVICOS DE SEGURANCA/VIGILANCIA - RUDDER EQUIPAMENTOS SISTEMAS SEGURANCA LTDADT COMPET: 10/02/2021'
$ws2.Range("D146").Value = 'ADM-01758/2020 - 11/02/2021-130974548-19 - SERV.FORNEC ENERGIA ELETRICA-CONTRATO CEEE 16 - ENERGIA ELETRICA - CEEE-CIA ESTADUAL DE ENERGIA ELETRICA DO RSDT COMPET: 11/02/2021'
$ws2.Range("D147").Value = 'ADM-01761/2020 - 17/02/2021-8139-17 - SERV. POSTAIS-CONTRATO FILIZOLA 18 - DESPESAS POSTAGEM E CARTORIO - FILIZOLA CENTRO DE SERVICOS LTDADT COMPET: 17/02/2021'
$ws2.Range("D148").Value = 'ADM-01763/2020 - PAPELARIA BRASIL LTDA NFC-E 583588-19 - MATERIAL DE ESCRITORIO - DT COMPET: 11/02/2021'
$ws2.Range("D149").Value = 'VLR PARCELA 07/12 DA AQUISIÇÃO DE 19 LICENÇAS ANTIVIRUS WF TECNOLOGIA LTDA'
$ws2.Range("D150").Value = 'VLR PARCELA 07/12 DA RENOVAÇÃO ANUAL DE UMA LICENÇA ANTIVIRUS WF TECNOLOGIA LTDA'
$ws2.Range("D151").Value = 'ADM-00001/2021 - ELABORACAO DA FOLHA DE PAGAMENTO E OBRIGACOES SOCIAIS 60, CFE NFSE N  2021-91 - FOLHA PAGAMENTO E CONSULT EM RH - MORESCO CONTABILIDADE S/SDT COMPET: 26/02/2021'
$ws2.Range("D152").Value = 'ADM-00004/2021 - DESP DESLOCAMENTO AO BANRISUL RETIRADA CARTAO ASSINATURAS DIRETORIA - CONDUCAO, TRANSPORTE E TELEENTREGA - LUANA RODRIGUES DE OLIVEIRADT COMPET: 25/02/2021'
$ws2.Range("D153").Value = 'ADM-00008/2021 - 02-2021 - PGTO PENSAO ALIMENTICIA - ELISANGELA AGUIRRE DOS SANTOSDT COMPET: 26/02/2021'
$ws2.Range("D154").Value = 'ADM-00008/2021 - 02-2021 - PGTO SALARIOS - ALEXANDRE VILLELA ILHADT COMPET: 26/02/2021'
$ws2.Range("D155").Value = 'ADM-00008/2021 - 02-2021 - PGTO SALARIOS - CASSIO ZARPELONDT COMPET: 26/02/2021'
$ws2.Range("D156").Value = 'ADM-00008/2021 - 02-2021 - PGTO SALARIOS - CATIA BERGER ROLIMDT COMPET: 26/02/2021'
$ws2.Range("D157").Value = 'ADM-00008/2021 - 02-2021 - PGTO SALARIOS - FLAVIO CARDOSO PINTO DA ROCHADT COMPET: 26/02/2021'
$ws2.Range("D158").Value = 'ADM-00008/2021 - 02-2021 - PGTO SALARIOS - LIGIAMAR FRIZZODT COMPET: 26/02/2021'
$ws2.Range("D159").Value = 'ADM-00008/2021 - 02-2021 - PGTO SALARIOS - MARCIA AMORIM MOREIRADT COMPET: 26/02/2021'
$ws2.Range("D160").Value = 'ADM-00008/2021 - 02-2021 - PGTO SALARIOS - NEIVA MINUSSI BIDINOTTODT COMPET: 26/02/2021'
$ws2.Range("D161").Value = 'ADM-00008/2021 - 02-2021 - PGTO SALARIOS - NICIA LOUYNEE MOREIRA WELLAUSEN PINTODT COMPET: 26/02/2021'
$ws2.Range("D162").Value = 'ADM-00008/2021 - 02-2021 - PGTO SALARIOS - SABRINA GIACOMONI COMELLIDT COMPET: 26/02/2021'
$ws2.Range("D163").Value = 'ADM-00009/2021 - 02-2021 - PGTO GRATIFICACAO DE REPRESENTACAO - ALVARO ROQUE KERN JUNQUEIRADT COMPET: 26/02/2021'
$ws2.Range("D164").Value = 'ADM-00009/2021 - 02-2021 - PGTO GRATIFICACAO DE REPRESENTACAO - JOSE PEDRO OSORIO MENDINADT COMPET: 26/02/2021'
$ws2.Range("D165").Value = 'ADM-00009/2021 - 02-2021 - PGTO GRATIFICACAO DE REPRESENTACAO - RICARDO ALTAIR SCHWARZDT COMPET: 26/02/2021'
$ws2.Range("D166").Value = 'ADM-00013/2021 - PG 05 DIAS DE FERIAS - FERIAS A PAGAR - LIGIAMAR FRIZZODT COMPET: 18/02/2021'
$ws2.Range("D176").Value = 'ADM-00015/2021 - SEGURO EMPRESARIAL PARA EMPREGADOS 51 - SEGURO FUNCIONARIOS - SOMPO SEGUROS S/ADT COMPET: 19/02/2021'
$ws2.Range("D177").Value = 'ADM-00017/2021 - DESLOCAMENTO DE TAXI RESIDENCIA ATE SEDE FAPERS-RECIBO EPTC N.448430 - CONDUCAO, TRANSPORTE E TELEENTREGA - ALEXANDRE VILLELA ILHADT COMPET: 24/02/2021'
$ws2.Range("D178").Value = 'ADM-00018/2021 - PG MAGAZINE LUIZA - NFC-E 516610 E 031958 - CARGA CHIC NUCLEO ADM FINANC PARA EMPRESTIMOS - SERVICOS DE TELEFONIA - DT COMPET: 24/02/2021'
$ws2.Range("D179").Value = 'ADM-00019/2021 - PG CLARO RS - CARGA CHIC GESTAO DOS PLANOS - SERVICOS DE TELEFONIA - DT COMPET: 24/02/2021'
$ws2.Range("D180").Value = 'ADM-00056/2021 - CONCESSAO DO PAT-VALE RANCHO 40 - VALE ALIMENTACAO - PATRONAL - TICKET SERVICOS S.A.DT COMPET: 26/02/2021'
$ws2.Range("D181").Value = 'ADM-00058/2021 - CONCESSAO PAT- VALE REFEICAO 41 - VALE REFEICAO-PATRONAL - TICKET SERVICOS S.A.DT COMPET: 26/02/2021'
$ws2.Range("D182").Value = 'ADM-00060/2021 - CONCESSAO TRI - VALE TRASPORTE 42 - VALE TRANSPORTE - ASSOC EMPRESAS TRANSPORTES PASSAGEIROS DE POADT COMPET: 26/02/2021'
$ws2.Range("D183").Value = 'ADM-00062/2021 - PG 4º TABELIONATO DE NOTAS - REC 1729570- DESPESAS DE AUTENTICACOES - DT COMPET: 25/02/2021'
$ws2.Range("D184").Value = 'ADM-01553/2020 - AQUISICAO DE MATERIAL DE HIGIENE E LIMPEZA - MATERIAL DE LIMPEZA - FABESUL COMERCIO DE SUPRIMENTOS LTDADT COMPET: 10/02/2021'
$ws2.Range("D185").Value = 'ADM-01553/2020 - AQUISICAO MATERIAL DE COPA E COZINHA - MATERIAL DE COPA/COZINHA - FABESUL COMERCIO DE SUPRIMENTOS LTDADT COMPET: 10/02/2021'
$ws2.Range("D186").Value = 'ADM-01553/2020 - AQUISICAO MATERIAL DE ESCRITORIO - MATERIAL DE ESCRITORIO/EXPEDIENTE - FABESUL COMERCIO DE SUPRIMENTOS LTDADT COMPET: 10/02/2021'
$ws2.Range("D187").Value = 'ADM-01554/2020 - AQUISICAO PAPEL TOALHA - MATERIAL DE LIMPEZA - RS PIRES COMERCIAL EIRELIDT COMPET: 12/02/2021'
$ws2.Range("D188").Value = 'ADM-01556/2020 - SERV. E-MAIL MARKTING-CONTRATO KINGHOST 14 NFSE.2021/12865 - SERVICOS DE INFORMATICA/INTERNET - CYBERWEB NETWORKS LTDADT COMPET: 05/02/2021'
$ws2.Range("D189").Value = 'ADM-01560/2020 - SERV. TECNICO SENSOR ALARME-CONTRATO RUDDER 11 - SERVICOS E MAT DE SEGURANCA/VIGILANCIA - RUDDER EQUIPAMENTOS SISTEMAS SEGURANCA LTDADT COMPET: 12/02/2021'
$ws2.Range("D190").Value = 'ADM-01562/2020 - SERVICOS DE ARMAZENAGEM/DIGITALIZACAO DE DOCUMENTOS 57 - GUARDA CONSERVACAO DE DOCUMENTOS - METROFILE BRASIL GESTAO DA INFORMACAO LTDADT COMPET: 05/02/2021'
$ws2.Range("D191").Value = 'ADM-01570/2020 - SERV VIVO FIXO - CONTRATO VIVO TELEFONICA 26 - SERVICOS DE TELEFONIA - TELEFONICA BRASIL S/A - VIVODT COMPET: 05/02/2021'
$ws2.Range("D192").Value = 'ADM-01574/2020 - SERV. MONITORAMENTO-TELEALARME-CONTRATO RUDDER 11 - SERVICOS DE SEGURANCA/VIGILANCIA - RUDDER EQUIPAMENTOS SISTEMAS SEGURANCA LTDADT COMPET: 12/02/2021'
$ws2.Range("D193").Value = 'ADM-01576/2020 - SERV. MONITORAMENTO-TELEALARME-CONTRATO RUDDER 11 - SERVICOS DE SEGURANCA/VIGILANCIA - RUDDER EQUIPAMENTOS SISTEMAS SEGURANCA LTDADT COMPET: 12/02/2021'
$ws2.Range("D194").Value = 'ADM-01587/2020 - SERV. TECNICO SENSOR ALARME-CONTRATO RUDDER 11 - SERVICOS E MAT DE SEGURANCA/VIGILANCIA - RUDDER EQUIPAMENTOS SISTEMAS SEGURANCA LTDADT COMPET: 24/02/2021'
$ws2.Range("D195").Value = 'ADM-01596/2020 - SERV. AGUA E ESGOTO-CONTRATO DMAE 17 - DMAE - DEPARTAMENTO MUNICIPAL DE AGUA E ESGOTOSDT COMPET: 10/02/2021'
$ws2.Range("D196").Value = 'ADM-01600/2020 - PARCELA 2 DE 4 SEGURO RESPONSAL CIVIL-D&O 56 - SEGURO RESPONSAB. CIVIL - AIG SEGUROS BRASIL S.A.DT COMPET: 26/02/2021'
$ws2.Range("D197").Value = 'ADM-01606/2020 - SERV NET - CONTRATO VIVO TELEFONICA 22 - SERVICOS DE TELEFONIA - TELEFONICA BRASIL S/A - VIVODT COMPET: 05/02/2021'
$ws2.Range("D198").Value = 'ADM-01608/2020 - SERV VIVO 0800 - CONTRATO VIVO TELEFONICA 23 - SERVICOS DE TELEFONIA - TELEFONICA BRASIL S/A - VIVODT COMPET: 05/02/2021'
$ws2.Range("D199").Value = 'ADM-01610/2020 - SERV. LOCACAO IMPRESSORAS/MULTIFUNC-CONTRATO PRINTMAX 12 - SERVICOS DE INFORMATICA/INTERNET - PRINTMAX GERENCIAMENTO DE DOCUMENTOS EIRELIDT COMPET: 05/02/2021'
$ws2.Range("D200").Value = 'ADM-01612/2020 - SERV TRANSPORTE PASSAG - CONTRATO CABIFY 26 - SERVICOS DE TRANSPORTE PASSAGEIROS - CABIFY AGENCIA DE SERVICOS DE TRANSPORTE DE PASSAGEIROS LTDADT COMPET: 05/02/2021'
$ws2.Range("D201").Value = 'ADM-01614/2020 - SERV. ADVOC. E DE CONSULTORIA JURIDICA-CONTRATO BARCELLOS 2 - CONSULTORIA JURIDICA - BARCELLOS ADVOCACIA EMPRESARIALDT COMPET: 12/02/2021'
$ws2.Range("D202").Value = 'ADM-01616/2020 - SERV PREST- REEMB DESPESAS POSTAGEM E CARTORIO 2 - DESPESAS POSTAGEM E CARTORIO - BARCELLOS ADVOCACIA EMPRESARIALDT COMPET: 12/02/2021'
$ws2.Range("D203").Value = 'ADM-01618/2020 - CONSULTORIA FINANCEIRA E DE INVESTIMENTOS 59 - CONSULTORIA INVESTIMENTOS - C R NEUENFELDT MEDT COMPET: 05/02/2021'
$ws2.Range("D204").Value = 'ADM-01620/2020 - ENTREGA DE MATERIAL ATRAVES DE MOOBOY - CONDUCAO, TRANSPORTE E TELEENTREGA - MAXIMO MOTOBOY EXPRESSO LTDADT COMPET: 05/02/2021'
$ws2.Range("D205").Value = 'ADM-01621/2020 - DESP.AUTENTICACOES 4.TABELIONATO REC. N.1726460 - DESPESAS POSTAGEM E CARTORIO - LUANA RODRIGUES DE OLIVEIRADT COMPET: 05/02/2021'
$ws2.Range("D206").Value = 'ADM-01621/2020 - DESPESAS DESLOCAMENTO UBER - CONDUCAO, TRANSPORTE E TELEENTREGA - LUANA RODRIGUES DE OLIVEIRADT COMPET: 05/02/2021'
$ws2.Range("D207").Value = 'ADM-01628/2020 - 01-2021 - FGTS - CAIXA ECON. FEDERALDT COMPET: 05/02/2021'
$ws2.Range("D208").Value = 'ADM-01629/2020 - 01-2021 - IRRF S/ TRABALHO ASSALARIADO - RECEITA FEDERAL DO BRASILDT COMPET: 19/02/2021'
$ws2.Range("D209").Value = 'ADM-01630/2020 - 01-2021 - INSS - INSS - INSTITUTO NACIONAL DO SEGURO SOCIALDT COMPET: 19/02/2021'
$ws2.Range("D210").Value = 'ADM-01635/2020 - CONVENIO AGENTE DE INTEGRACAO BOLSA DE ESTAGIO 54 - PGTO BOLSA ESTAGIO - CIEERS - CENTRO DE INTEGRACAO EMPRESA ESCOLADT COMPET: 01/02/2021'
$ws2.Range("D211").Value = 'ADM-01635/2020 - CONVENIO AGENTE DE INTEGRACAO BOLSA DE ESTAGIO 54 - PGTO TX ADM BOLSA ESTAGIO - CIEERS - CENTRO DE INTEGRACAO EMPRESA ESCOLADT COMPET: 01/02/2021'
$ws2.Range("D212").Value = 'ADM-01640/2020 - LICENCA DE USO DE SOFTWARE-CONTRATO INTECH 1 - ASSES/CONSULT INFORMATICA - INTECH SOLUCOES EM TECNOLOGIA DA INFORMACAO LTDADT COMPET: 10/02/2021'
$ws2.Range("D213").Value = 'ADM-01642/2020 - CONSULT, MANUT E SERV TECNICOS EM INFORMTICA 15 - SERVICOS DE INFORMATICA/INTERNET - WF TECNOLOGIA LTDADT COMPET: 12/02/2021'
$ws2.Range("D214").Value = 'ADM-01644/2020 - SERV. POSTAIS-CONTRATO FILIZOLA 18 - DESPESAS POSTAGEM E CARTORIO - FILIZOLA CENTRO DE SERVICOS LTDADT COMPET: 26/02/2021'
$ws2.Range("D215").Value = 'ADM-01646/2020 - LICENCA SOFTWARE QUANTUM INDICADORES INVESTIMENTOS 37 - CONSULTORIA INVESTIMENTOS - PARAMITA TECNOLOGIA CONSULTORIA FINANCEIRA LTDADT COMPET: 10/02/2021'
$ws2.Range("D216").Value = 'ADM-01648/2020 - SERVICOS DE LIMPEZA PREDIAL-CONTRATO SULCLEAN 5 - SERVICOS DE LIMPEZA - SULCLEAN SERVICOS LTDADT COMPET: 10/02/2021'
$ws2.Range("D217").Value = 'ADM-01650/2020 - CONTRIB ASSOCIATIVA - CONTRATO ABRAPP 20 - ENTIDADES E ASSOCIACOES - ASSOCIACAO BRAS. DE ENTIDADES FECHADAS PREVIDENCIA PRIVADADT COMPET: 10/02/2021'
$ws2.Range("D218").Value = 'ADM-01652/2020 - PREST SERV MAPEAMENTO,OTIMIZ PROCESSOS 50 - ASSES/CONSULT MAPEAM,OTIMIZ PROCESSOS - MENTHOR CONSULTORIA E CAPACITACAO EMPRESARIAL LTDADT COMPET: 10/02/2021'
$ws2.Range("D219").Value = 'ADM-01674/2020 - CONSULTORIA ATUARIAL-METODO ATUARIAL 58 - CONSULTORIA ATUARIAL - METODO ATUARIAL SOCIEDADE SIMPLES LTDADT COMPET: 05/02/2021'
$ws2.Range("D220").Value = 'ADM-01680/2020 - ASSES. E CONSUL. DE INVESTIMENTOS-CONTRATO MARCO & MARCO 6 - CONSULTORIA INVESTIMENTOS - MARCO & MARCO CONSULTORES FINANCEIROS ASSOCIADOS LTDADT COMPET: 10/02/2021'
$ws2.Range("D221").Value = 'ADM-01687/2020 - SERV. HOSPEDAGEM DE SITES-CONTRATO KINGHOST 14 - SERVICOS DE INFORMATICA/INTERNET - CYBERWEB NETWORKS LTDADT COMPET: 05/02/2021'
$ws2.Range("D222").Value = 'ADM-01689/2020 - SERV VIVO INTERNET- CONTRATO VIVO TELEFONICA 24 - SERVICOS DE TELEFONIA - TELEFONICA BRASIL S/A - VIVODT COMPET: 12/02/2021'
$ws2.Range("D223").Value = 'ADM-01695/2020 - ASSESS/CONSULT - SERV AREA COMUNICACAO/REDE SOCIAIS 33 - ASSES/CONSULT GESTAO/PLANEJ. ESTRATEGICO - MIRADOR ASSESSORIA ATUARIAL LTDADT COMPET: 10/02/2021'
$ws2.Range("D224").Value = 'ADM-01697/2020 - ASSES/CONSULT-SUPORTE TECNICO SITE FAPERS 34 - SERVICOS TECNICO SUPORTE SITE - MIRADOR ASSESSORIA ATUARIAL LTDADT COMPET: 10/02/2021'
$ws2.Range("D225").Value = 'ADM-01699/2020 - TRANSPORTE DE UBER RESIDENCIA ATE FAPERS - CONDUCAO, TRANSPORTE E TELEENTREGA - CASSIO ZARPELONDT COMPET: 10/02/2021'
$ws2.Range("D226").Value = 'ADM-01701/2020 - PGTO EMOLUMENTOS NO CARTORIO DE REGISTRO ESPECIAL- DIFERENCA - DESPESAS POSTAGEM E CARTORIO - CATIA BERGER ROLIMDT COMPET: 10/02/2021'
$ws2.Range("D227").Value = 'ADM-01702/2020 - DELOCAMENTO DE ATA VIA RODOVIARIA - CONDUCAO, TRANSPORTE E TELEENTREGA - DIRLEI MATOS DE SOUZADT COMPET: 10/02/2021'
$ws2.Range("D228").Value = 'ADM-01703/2020 - PREST SERVICOS MEDICOS-CO-PARTICIPACAO 36 - PLANO DE SAUDE - EMPREGADOS - UNIMED PORTO ALEGRE SOCIEDADE COOPERATIVA DE TRABALHO MEDICODT COMPET: 10/02/2021'
$ws2.Range("D229").Value = 'ADM-01705/2020 - PREST SERVICOS DE CONSULTAS MEDICAS 35 - PLANO DE SAUDE - EMPREGADOS - UNIMED PORTO ALEGRE SOCIEDADE COOPERATIVA DE TRABALHO MEDICODT COMPET: 10/02/2021'
$ws2.Range("D230").Value = 'ADM-01705/2020 - PREST SERVICOS DE CONSULTAS MEDICAS 35 - PLANO DE SAUDE - PATRONAL - UNIMED PORTO ALEGRE SOCIEDADE COOPERATIVA DE TRABALHO MEDICODT COMPET: 10/02/2021'
$ws2.Range("D231").Value = 'ADM-01714/2020 - PREST SERV ODONTOLOGICOS 34 - PLANO ODONTOLOGICO - EMPREGADOS - UNIODONTO PORTO ALEGRE - COOPERATIVA ODONTOLOGICA LTDADT COMPET: 10/02/2021'
$ws2.Range("D232").Value = 'ADM-01714/2020 - PREST SERV ODONTOLOGICOS 34 - PLANO ODONTOLOGICO - PATRONAL - UNIODONTO PORTO ALEGRE - COOPERATIVA ODONTOLOGICA LTDADT COMPET: 10/02/2021'
$ws2.Range("D233").Value = 'ADM-01717/2020 - PAGAMENTO DE DESPESAS DE CARTORIO - DESPESAS POSTAGEM E CARTORIO - CATIA BERGER ROLIMDT COMPET: 10/02/2021'
$ws2.Range("D234").Value = 'ADM-01717/2020 - PGTO TAXIS E APLICATIVOS PARA DESCOLAMENTO AO CARTORIO - CONDUCAO, TRANSPORTE E TELEENTREGA - CATIA BERGER ROLIMDT COMPET: 10/02/2021'
$ws2.Range("D235").Value = 'ADM-01718/2020 - PG 4º TABELIONATO DE NOTAS - DESPESAS DE XEROX E AUTENTICACOES CFE.RECIBO N  1727468 - DESPESAS POSTAGEM E CARTORIO - DT COMPET: 10/02/2021'
$ws2.Range("D236").Value = 'ADM-01720/2020 - SERV MOVEL CELULARES - CONTRATO VIVO TELEFONICA 25 - SERVICOS DE TELEFONIA - TELEFONICA BRASIL S/A - VIVODT COMPET: 24/02/2021'
$ws2.Range("D237").Value = 'ADM-01760/2020 - IMPRESSAO DE CARTAO DE VISITAS CFE.REC.N  01342766 - GRAFICA E PUBLICIDADE - GRAFICA ANSDT COMPET: 12/02/2021'
$ws2.Range("D238").Value = 'ADM-01761/2020 - SERV. POSTAIS-CONTRATO FILIZOLA 18 - DESPESAS POSTAGEM E CARTORIO - FILIZOLA CENTRO DE SERVICOS LTDADT COMPET: 24/02/2021'
$ws2.Range("D239").Value = 'ADM-01763/2020 - PG PAPELARIA BRASIL - NFC-E 583588 - MATERIAL DE ESCRITORIO - DT COMPET: 18/02/2021'
$ws2.Range("D240").Value = 'CONT-00001/2021 - CONSIGNATARIO - ISSQN - RECEITA FEDERAL DO BRASILDT COMPET: 10/02/2021'
$ws2.Range("D241").Value = 'CONT-00002/2021 - CONSIGNATARIO - PIS/COFINS/CSLL SERV PRESTADO - RECEITA FEDERAL DO BRASILDT COMPET: 19/02/2021'
$ws2.Range("D242").Value = 'CONT-00003/2021 - CONSIGNATARIO - IRRF S/ SERV PRESTADO - RECEITA FEDERAL DO BRASILDT COMPET: 19/02/2021'
$ws2.Range("D243").Value = 'CONT-00005/2021 - CONSIGNATARIO - INSS S/ SERVICO - INSS - INSTITUTO NACIONAL DO SEGURO SOCIALDT COMPET: 19/02/2021'
$ws2.Range("D244").Value = 'VLR CSLL/PIS/COFINS S/SERVIÇO SULCLEAN SERV LTDA REF JANEIRO/21'
$ws2.Range("D245").Value = 'INV-00004/2021 - FMULT-0148/2019  VINCI VALOREM FIM -  RF - BOLETA:/ - FMULT FUNDO MULTIMERCADO - RESGATE - VINCI VALOREM FUNDO DE INVESTIMENTO MULTIMERCADODT COMPET: 03/02/2021'
$ws2.Range("D246").Value = 'INV-00018/2021 - FMULT-0047/2020  VINCI VALOREM FIM -  RF - BOLETA:/ - FMULT FUNDO MULTIMERCADO - RESGATE - VINCI VALOREM FUNDO DE INVESTIMENTO MULTIMERCADODT COMPET: 17/02/2021'
$ws2.Range("D247").Value = 'INV-00036/2021 - NTN-B-0043/2013  STNC -  RF - BOLETA:730/2012 - NTN-B - RESGATE - STNCDT COMPET: 17/02/2021'
$ws2.Range("D248").Value = 'VLR CSLL/PIS/COFINS S/SERVIÇO SULCLEAN SERV LTDA REF JANEIRO/21'
$ws2.Range("D249").Value = 'VLR CSLL/PIS/COFINS SERVIÇO METODO ATUARIAL LTDA REF JANEIRO/2020'
$ws2.Range("D250").Value = 'VLR CSLL/PIS/COFINS SERVIÇO METODO ATUARIAL LTDA REF JANEIRO/2021'
$ws2.Range("D251").Value = 'VLR CSLL/PIS/COFINS BARCELLOS ADVOC EMPRES REF JANEIRO/21'
$ws2.Range("D252").Value = 'VLR CSLL/PIS/COFINS S/SERVIÇO DE INTECH SOLUÇÃO EM TECNO DA INFO LTDA REF JANEIRO/21'
$ws2.Range("D253").Value = 'RET CSLL/PIS/COFINS S/SERVIÇO DE MIRADOR ASSESSORIA ATUARIAL LTDA REF JANEIRO/21'
$ws2.Range("D254").Value = 'AMORTIZAÇÃO - NATUREZA: SOFTWARE - INTANGIVEL - TOMBAMENTO: 761 - 1'
$ws2.Range("D255").Value = 'AMORTIZAÇÃO - NATUREZA: SOFTWARE - INTANGIVEL - TOMBAMENTO: 767 - 1'
$ws2.Range("D256").Value = 'AMORTIZAÇÃO - NATUREZA: SOFTWARE - INTANGIVEL - TOMBAMENTO: 768 - 1'
$ws2.Range("D257").Value = 'AMORTIZAÇÃO - NATUREZA: SOFTWARE - INTANGIVEL - TOMBAMENTO: 769 - 1'
$ws2.Range("D258").Value = 'AMORTIZAÇÃO - NATUREZA: SOFTWARE - INTANGIVEL - TOMBAMENTO: 770 - 1'
$ws2.Range("D259").Value = 'AMORTIZAÇÃO - NATUREZA: SOFTWARE - INTANGIVEL - TOMBAMENTO: 771 - 1'
$ws2.Range("D260").Value = 'AMORTIZAÇÃO - NATUREZA: SOFTWARE - INTANGIVEL - TOMBAMENTO: 772 - 1'
$ws2.Range("D261").Value = 'AMORTIZAÇÃO - NATUREZA: SOFTWARE - INTANGIVEL - TOMBAMENTO: 773 - 1'
$ws2.Range("D262").Value = 'AMORTIZAÇÃO - NATUREZA: SOFTWARE - INTANGIVEL - TOMBAMENTO: 774 - 1'
$ws2.Range("D263").Value = 'AMORTIZAÇÃO - NATUREZA: SOFTWARE - INTANGIVEL - TOMBAMENTO: 775 - 1'
$ws2.Range("D264").Value = 'AMORTIZAÇÃO - NATUREZA: SOFTWARE - INTANGIVEL - TOMBAMENTO: 776 - 1'
$ws2.Range("D265").Value = 'AMORTIZAÇÃO - NATUREZA: SOFTWARE - INTANGIVEL - TOMBAMENTO: 812 - 1'
$ws2.Range("D266").Value = 'DEPRECIAÇÃO - NATUREZA: COMPUTADORES E PERIFÉRICOS - TOMBAMENTO: 760 - 1'
$ws2.Range("D267").Value = 'DEPRECIAÇÃO - NATUREZA: COMPUTADORES E PERIFÉRICOS - TOMBAMENTO: 762 - 1'
$ws2.Range("D268").Value = 'DEPRECIAÇÃO - NATUREZA: COMPUTADORES E PERIFÉRICOS - TOMBAMENTO: 764 - 1'
$ws2.Range("D269").Value = 'DEPRECIAÇÃO - NATUREZA: COMPUTADORES E PERIFÉRICOS - TOMBAMENTO: 777 - 1'
$ws2.Range("D270").Value = 'DEPRECIAÇÃO - NATUREZA: COMPUTADORES E PERIFÉRICOS - TOMBAMENTO: 778 - 1'
$ws2.Range("D271").Value = 'DEPRECIAÇÃO - NATUREZA: COMPUTADORES E PERIFÉRICOS - TOMBAMENTO: 779 - 1'
$ws2.Range("D272").Value = 'DEPRECIAÇÃO - NATUREZA: COMPUTADORES E PERIFÉRICOS - TOMBAMENTO: 781 - 1'
$ws2.Range("D273").Value = 'DEPRECIAÇÃO - NATUREZA: COMPUTADORES E PERIFÉRICOS - TOMBAMENTO: 782 - 1'
$ws2.Range("D274").Value = 'DEPRECIAÇÃO - NATUREZA: COMPUTADORES E PERIFÉRICOS - TOMBAMENTO: 783 - 1'
$ws2.Range("D275").Value = 'DEPRECIAÇÃO - NATUREZA: COMPUTADORES E PERIFÉRICOS - TOMBAMENTO: 784 - 1'
$ws2.Range("D276").Value = 'DEPRECIAÇÃO - NATUREZA: COMPUTADORES E PERIFÉRICOS - TOMBAMENTO: 787 - 1'
$ws2.Range("D277").Value = 'DEPRECIAÇÃO - NATUREZA: COMPUTADORES E PERIFÉRICOS - TOMBAMENTO: 788 - 1'
$ws2.Range("D278").Value = 'DEPRECIAÇÃO - NATUREZA: COMPUTADORES E PERIFÉRICOS - TOMBAMENTO: 789 - 1'
$ws2.Range("D279").Value = 'DEPRECIAÇÃO - NATUREZA: COMPUTADORES E PERIFÉRICOS - TOMBAMENTO: 806 - 1'
$ws2.Range("D280").Value = 'DEPRECIAÇÃO - NATUREZA: COMPUTADORES E PERIFÉRICOS - TOMBAMENTO: 807 - 1'
$ws2.Range("D281").Value = 'DEPRECIAÇÃO - NATUREZA: COMPUTADORES E PERIFÉRICOS - TOMBAMENTO: 808 - 1'
$ws2.Range("D282").Value = 'DEPRECIAÇÃO - NATUREZA: COMPUTADORES E PERIFÉRICOS - TOMBAMENTO: 809 - 1'
$ws2.Range("D283").Value = 'DEPRECIAÇÃO - NATUREZA: COMPUTADORES E PERIFÉRICOS - TOMBAMENTO: 810 - 1'
$ws2.Range("D284").Value = 'DEPRECIAÇÃO - NATUREZA: COMPUTADORES E PERIFÉRICOS - TOMBAMENTO: 811 - 1'
$ws2.Range("D285").Value = 'DEPRECIAÇÃO - NATUREZA: COMPUTADORES E PERIFÉRICOS - TOMBAMENTO: 814 - 1'
$ws2.Range("D286").Value = 'DEPRECIAÇÃO - NATUREZA: COMPUTADORES E PERIFÉRICOS - TOMBAMENTO: 817 - 1'
$ws2.Range("D287").Value = 'DEPRECIAÇÃO - NATUREZA: COMPUTADORES E PERIFÉRICOS - TOMBAMENTO: 818 - 1'
$ws2.Range("D288").Value = 'DEPRECIAÇÃO - NATUREZA: COMPUTADORES E PERIFÉRICOS - TOMBAMENTO: 819 - 1'
$ws2.Range("D289").Value = 'DEPRECIAÇÃO - NATUREZA: COMPUTADORES E PERIFÉRICOS - TOMBAMENTO: 820 - 1'
$ws2.Range("D290").Value = 'DEPRECIAÇÃO - NATUREZA: COMPUTADORES E PERIFÉRICOS - TOMBAMENTO: 821 - 1'
$ws2.Range("D291").Value = 'DEPRECIAÇÃO - NATUREZA: COMPUTADORES E PERIFÉRICOS - TOMBAMENTO: 824 - 1'
$ws2.Range("D292").Value = 'DEPRECIAÇÃO - NATUREZA: COMPUTADORES E PERIFÉRICOS - TOMBAMENTO: 825 - 1'
$ws2.Range("D293").Value = 'DEPRECIAÇÃO - NATUREZA: COMPUTADORES E PERIFÉRICOS - TOMBAMENTO: 826 - 1'
$ws2.Range("D294").Value = 'DEPRECIAÇÃO - NATUREZA: COMPUTADORES E PERIFÉRICOS - TOMBAMENTO: 827 - 1'
$ws2.Range("D295").Value = 'DEPRECIAÇÃO - NATUREZA: MAQUINAS E EQUIPAMENTOS - TOMBAMENTO: 751 - 1'
$ws2.Range("D296").Value = 'DEPRECIAÇÃO - NATUREZA: MAQUINAS E EQUIPAMENTOS - TOMBAMENTO: 756 - 1'
$ws2.Range("D297").Value = 'DEPRECIAÇÃO - NATUREZA: MAQUINAS E EQUIPAMENTOS - TOMBAMENTO: 763 - 1'
$ws2.Range("D298").Value = 'DEPRECIAÇÃO - NATUREZA: MAQUINAS E EQUIPAMENTOS - TOMBAMENTO: 822 - 1'
$ws2.Range("D299").Value = 'DEPRECIAÇÃO - NATUREZA: MOVEIS E UTENSILIOS - TOMBAMENTO: 752 - 1'
$ws2.Range("D300").Value = 'DEPRECIAÇÃO - NATUREZA: MOVEIS E UTENSILIOS - TOMBAMENTO: 753 - 1'
$ws2.Range("D301").Value = 'DEPRECIAÇÃO - NATUREZA: MOVEIS E UTENSILIOS - TOMBAMENTO: 754 - 1'
$ws2.Range("D302").Value = 'DEPRECIAÇÃO - NATUREZA: MOVEIS E UTENSILIOS - TOMBAMENTO: 765 - 1'
$ws2.Range("D303").Value = 'DEPRECIAÇÃO - NATUREZA: MOVEIS E UTENSILIOS - TOMBAMENTO: 766 - 1'
$ws2.Range("D304").Value = 'DEPRECIAÇÃO - NATUREZA: MOVEIS E UTENSILIOS - TOMBAMENTO: 790 - 1'
$ws2.Range("D305").Value = 'DEPRECIAÇÃO - NATUREZA: MOVEIS E UTENSILIOS - TOMBAMENTO: 791 - 1'
$ws2.Range("D306").Value = 'DEPRECIAÇÃO - NATUREZA: MOVEIS E UTENSILIOS - TOMBAMENTO: 792 - 1'
$ws2.Range("D307").Value = 'DEPRECIAÇÃO - NATUREZA: MOVEIS E UTENSILIOS - TOMBAMENTO: 793 - 1'
$ws2.Range("D308").Value = 'DEPRECIAÇÃO - NATUREZA: MOVEIS E UTENSILIOS - TOMBAMENTO: 794 - 1'
$ws2.Range("D309").Value = 'DEPRECIAÇÃO - NATUREZA: MOVEIS E UTENSILIOS - TOMBAMENTO: 795 - 1'
$ws2.Range("D310").Value = 'DEPRECIAÇÃO - NATUREZA: MOVEIS E UTENSILIOS - TOMBAMENTO: 796 - 1'
$ws2.Range("D311").Value = 'DEPRECIAÇÃO - NATUREZA: MOVEIS E UTENSILIOS - TOMBAMENTO: 797 - 1'
$ws2.Range("D312").Value = 'DEPRECIAÇÃO - NATUREZA: MOVEIS E UTENSILIOS - TOMBAMENTO: 798 - 1'
$ws2.Range("D313").Value = 'DEPRECIAÇÃO - NATUREZA: MOVEIS E UTENSILIOS - TOMBAMENTO: 799 - 1'
$ws2.Range("D314").Value = 'DEPRECIAÇÃO - NATUREZA: MOVEIS E UTENSILIOS - TOMBAMENTO: 800 - 1'
$ws2.Range("D315").Value = 'DEPRECIAÇÃO - NATUREZA: MOVEIS E UTENSILIOS - TOMBAMENTO: 801 - 1'
$ws2.Range("D316").Value = 'DEPRECIAÇÃO - NATUREZA: MOVEIS E UTENSILIOS - TOMBAMENTO: 802 - 1'
$ws2.Range("D317").Value = 'DEPRECIAÇÃO - NATUREZA: MOVEIS E UTENSILIOS - TOMBAMENTO: 803 - 1'
$ws2.Range("D318").Value = 'DEPRECIAÇÃO - NATUREZA: MOVEIS E UTENSILIOS - TOMBAMENTO: 804 - 1'
$ws2.Range("D319").Value = 'DEPRECIAÇÃO - NATUREZA: MOVEIS E UTENSILIOS - TOMBAMENTO: 813 - 1'
$ws2.Range("D320").Value = 'PG TELEFONICA BRASIL SA PARCELA 16/24 CFE NFS-E 3734276'
$ws2.Range("D321").Value = 'VLR FONTE CUSTEIO ADM DOS INVESTIMENTOS PGS '
$ws2.Range("D322").Value = 'VLR FONTE CUSTEIO ADM DOS INVESTIMENTOS PLANO MISTO '
$ws2.Range("D323").Value = 'VLR FONTE CUSTEIO ADM DOS INVESTIMENTOS PREVER '
$ws2.Range("D324").Value = 'LIQ PARC ALUGUEL SEDE FAPERS REF PARTE PBD-I JAN/21'
$ws2.Range("D325").Value = 'LIQ PARC ALUGUEL SEDE FAPERS REF PARTE PGS JAN/21'
$ws2.Range("D326").Value = 'LIQ PARC ALUGUEL SEDE FAPERS REF PARTE PLANO MISTO JAN/21'
$ws2.Range("D327").Value = 'LIQ PARC ALUGUEL SEDE FAPERS REF PARTE PREVER JAN/21'
$ws2.Range("D328").Value = 'LIQ PARC CUSTEIO ADM INVEST PGS REF JAN/21'
$ws2.Range("D329").Value = 'LIQ PARC CUSTEIO ADM INVEST PLANO MISTO REF JAN/21'
$ws2.Range("D330").Value = 'LIQ PARC CUSTEIO ADM INVEST PREVER REF JAN/21'
$ws2.Range("D331").Value = 'APROPR CONTR ASSISITIDOS PGS PATROC FAPERS N/MES'
$ws2.Range("D332").Value = 'APROPR CONTR ATIVOS PGS PATROC FAPERS N/MES'
$ws2.Range("D333").Value = 'APROPR CONTRIB ATIVOS PLANO MISTO PATROC FAPERS REF N/MES'
$ws2.Range("D334").Value = 'VLR LIQ CCA DA PATROC FAPERS P/ PGS '
$ws2.Range("D335").Value = 'VLR LIQ CCA DA PATROC FAPERS P/ PREVER '
$ws2.Range("D336").Value = 'VLR LIQ CCA DA PATROC FAPERS REF ASSISTIDOS P/ PGS '
$ws2.Range("D337").Value = 'VLR LIQ CCA DA PATROC FAPERS REF ASSISTIDOS P/ PLANO MISTO '
$ws2.Range("D338").Value = 'VLR LIQ CCA PARTIC ATIVOS DA PATROC FAPERS P/ PGS '
$ws2.Range("D339").Value = 'VLR LIQ CCA PARTIC ATIVOS DA PATROC FAPERS P/ PREVER '
$ws2.Range("D340").Value = 'VLR LIQ CONTR DA PATROC FAPERS REF ATIVOS P/ PLANO MISTO '
$ws2.Range("D341").Value = 'VLR LIQ CONTR DA PATROC FAPERS REF ATIVOS P/ PLANO MISTO DE JAN/21'
$ws2.Range("D342").Value = 'VLR LIQ CONTR PARTIC ATIVOS DA PATROC FAPERS P/ PLANO MISTO '
$ws2.Range("D343").Value = 'VLR LIQ CONTR PARTIC ATIVOS DA PATROC FAPERS P/ PLANO MISTO REF JAN/21'
$ws2.Range("D344").Value = 'VLR APLICAÇÃO  FIRF - 10/2021'
$ws2.Range("D345").Value = 'VLR APLICAÇÃO  FIRF - 2/2021'
$ws2.Range("D346").Value = 'VLR APLICAÇÃO  FIRF - 9/2021'
$ws2.Range("D347").Value = 'VLR APROP VARIAÇÃO NEGATIVA  FIRF - 2/2021'
$ws2.Range("D348").Value = 'VLR APROP VARIAÇÃO NEGATIVA  FIRF - 9/2021'
$ws2.Range("D349").Value = 'VLR APROP VARIAÇÃO NEGATIVA  FMULT - 11/2021'
$ws2.Range("D350").Value = 'VLR APROP VARIAÇÃO NEGATIVA  FMULT - 12/2021'
$ws2.Range("D351").Value = 'VLR APROP VARIAÇÃO NEGATIVA  FMULT - 125/2020'
$ws2.Range("D352").Value = 'VLR APROP VARIAÇÃO NEGATIVA  FMULT - 126/2020'
$ws2.Range("D353").Value = 'VLR APROP VARIAÇÃO NEGATIVA  FMULT - 148/2019'
$ws2.Range("D354").Value = 'VLR APROP VARIAÇÃO NEGATIVA  FMULT - 148/2020'
$ws2.Range("D355").Value = 'VLR APROP VARIAÇÃO NEGATIVA  FMULT - 149/2020'
$ws2.Range("D356").Value = 'VLR APROP VARIAÇÃO NEGATIVA  FMULT - 158/2019'
$ws2.Range("D357").Value = 'VLR APROP VARIAÇÃO NEGATIVA  FMULT - 169/2020'
$ws2.Range("D358").Value = 'VLR APROP VARIAÇÃO NEGATIVA  FMULT - 171/2020'
$ws2.Range("D359").Value = 'VLR APROP VARIAÇÃO NEGATIVA  FMULT - 176/2019'
$ws2.Range("D360").Value = 'VLR APROP VARIAÇÃO NEGATIVA  FMULT - 47/2020'
$ws2.Range("D361").Value = 'VLR APROP VARIAÇÃO NEGATIVA  FMULT - 64/2020'
$ws2.Range("D362").Value = 'VLR APROP VARIAÇÃO NEGATIVA  FMULT - 84/2020'
$ws2.Range("D363").Value = 'VLR APROP VARIAÇÃO NEGATIVA  FMULT - 96/2020'
$ws2.Range("D364").Value = 'VLR APROP VARIAÇÃO NEGATIVA  FMULT - 97/2020'
$ws2.Range("D365").Value = 'VLR APROP VARIAÇÃO NEGATIVA  FMULT - 98/2020'
$ws2.Range("D366").Value = 'VLR APROP VARIAÇÃO NEGATIVA  NTN-B - 31/2013'
$ws2.Range("D367").Value = 'VLR APROP VARIAÇÃO NEGATIVA  NTN-B - 32/2013'
$ws2.Range("D368").Value = 'VLR APROP VARIAÇÃO NEGATIVA  NTN-B - 34/2013'
$ws2.Range("D369").Value = 'VLR APROP VARIAÇÃO NEGATIVA  NTN-B - 4/2013'
$ws2.Range("D370").Value = 'VLR APROP VARIAÇÃO NEGATIVA  NTN-B - 5/2013'
$ws2.Range("D371").Value = 'VLR APROP VARIAÇÃO NEGATIVA  NTN-B - 507/2013'
$ws2.Range("D372").Value = 'VLR APROP VARIAÇÃO NEGATIVA  NTN-B - 510/2013'
$ws2.Range("D373").Value = 'VLR APROP VARIAÇÃO NEGATIVA  NTN-B - 522/2013'
$ws2.Range("D374").Value = 'VLR APROP VARIAÇÃO NEGATIVA  NTN-B - 523/2013'
$ws2.Range("D375").Value = 'VLR APROP VARIAÇÃO NEGATIVA  NTN-B - 6/2013'
$ws2.Range("D376").Value = 'VLR APROP VARIAÇÃO POSITIVA  FIEST - 170/2020'
$ws2.Range("D377").Value = 'VLR APROP VARIAÇÃO POSITIVA  FIEST - 4/2020'
$ws2.Range("D378").Value = 'VLR APROP VARIAÇÃO POSITIVA  FIEST - 8/2021'
$ws2.Range("D379").Value = 'VLR APROP VARIAÇÃO POSITIVA  FIRF - 2/2021'
$ws2.Range("D380").Value = 'VLR APROP VARIAÇÃO POSITIVA  FIRF - 9/2021'
$ws2.Range("D381").Value = 'VLR APROP VARIAÇÃO POSITIVA  FMULT - 11/2021'
$ws2.Range("D382").Value = 'VLR APROP VARIAÇÃO POSITIVA  FMULT - 12/2021'
$ws2.Range("D383").Value = 'VLR APROP VARIAÇÃO POSITIVA  FMULT - 125/2020'
$ws2.Range("D384").Value = 'VLR APROP VARIAÇÃO POSITIVA  FMULT - 126/2020'
$ws2.Range("D385").Value = 'VLR APROP VARIAÇÃO POSITIVA  FMULT - 148/2019'
$ws2.Range("D386").Value = 'VLR APROP VARIAÇÃO POSITIVA  FMULT - 148/2020'
$ws2.Range("D387").Value = 'VLR APROP VARIAÇÃO POSITIVA  FMULT - 149/2020'
$ws2.Range("D388").Value = 'VLR APROP VARIAÇÃO POSITIVA  FMULT - 158/2019'
$ws2.Range("D389").Value = 'VLR APROP VARIAÇÃO POSITIVA  FMULT - 171/2020'
$ws2.Range("D390").Value = 'VLR APROP VARIAÇÃO POSITIVA  FMULT - 176/2019'
$ws2.Range("D391").Value = 'VLR APROP VARIAÇÃO POSITIVA  FMULT - 47/2020'
$ws2.Range("D392").Value = 'VLR APROP VARIAÇÃO POSITIVA  FMULT - 64/2020'
$ws2.Range("D393").Value = 'VLR APROP VARIAÇÃO POSITIVA  FMULT - 84/2020'
$ws2.Range("D394").Value = 'VLR APROP VARIAÇÃO POSITIVA  FMULT - 96/2020'
$ws2.Range("D395").Value = 'VLR APROP VARIAÇÃO POSITIVA  FMULT - 97/2020'
$ws2.Range("D396").Value = 'VLR APROP VARIAÇÃO POSITIVA  FMULT - 98/2020'
$ws2.Range("D397").Value = 'VLR APROP VARIAÇÃO POSITIVA  NTN-B - 31/2013'
$ws2.Range("D398").Value = 'VLR APROP VARIAÇÃO POSITIVA  NTN-B - 32/2013'
$ws2.Range("D399").Value = 'VLR APROP VARIAÇÃO POSITIVA  NTN-B - 34/2013'
$ws2.Range("D400").Value = 'VLR APROP VARIAÇÃO POSITIVA  NTN-B - 4/2013'
$ws2.Range("D401").Value = 'VLR APROP VARIAÇÃO POSITIVA  NTN-B - 5/2013'
$ws2.Range("D402").Value = 'VLR APROP VARIAÇÃO POSITIVA  NTN-B - 507/2013'
$ws2.Range("D403").Value = 'VLR APROP VARIAÇÃO POSITIVA  NTN-B - 510/2013'
$ws2.Range("D404").Value = 'VLR APROP VARIAÇÃO POSITIVA  NTN-B - 522/2013'
$ws2.Range("D405").Value = 'VLR APROP VARIAÇÃO POSITIVA  NTN-B - 523/2013'
$ws2.Range("D406").Value = 'VLR APROP VARIAÇÃO POSITIVA  NTN-B - 6/2013'
$ws2.Range("D407").Value = 'VLR REF RECLASSIFICAÇÃO FUNDO KINEA IPCA FIC FIM [KINEAABS] - DE MULTIMERCADO PARA RENDA FIXA'
$ws2.Range("D408").Value = 'VLR RESGATE A RECEBER  FMULT - 148/2019'
$ws2.Range("D409").Value = 'VLR RESGATE A RECEBER  FMULT - 158/2019'
$ws2.Range("D410").Value = 'VLR RESGATE A RECEBER  FMULT - 176/2019'
$ws2.Range("D411").Value = 'VLR RESGATE A RECEBER  FMULT - 47/2020'
$ws2.Range("D412").Value = 'VLR RESGATE A RECEBER  NTN-B - BOLETA Nº 723/2012'
$ws2.Range("D413").Value = 'VLR RESGATE A RECEBER  NTN-B - BOLETA Nº 724/2012'
$ws2.Range("D414").Value = 'VLR RESGATE A RECEBER  NTN-B - BOLETA Nº 725/2012'
$ws2.Range("D415").Value = 'VLR RESGATE A RECEBER  NTN-B - BOLETA Nº 726/2012'
$ws2.Range("D416").Value = 'VLR RESGATE A RECEBER  NTN-B - BOLETA Nº 727/2012'
$ws2.Range("D417").Value = 'VLR RESGATE A RECEBER  NTN-B - BOLETA Nº 728/2012'
$ws2.Range("D418").Value = 'VLR RESGATE A RECEBER  NTN-B - BOLETA Nº 729/2012'
$ws2.Range("D419").Value = 'RV-BOLETA O-24/2021 PAPEL GEO GLOBAL EQ FIA IE FND'
$ws2.Range("D420").Value = 'RV-BOLETA O-29/2021 PAPEL SAFRA CONS AMER FIA FND'
$ws2.Range("D421").Value = 'VALORES A LIQUIDAR RV-BOLETA O-29/2021'
$ws2.Range("D422").Value = 'VARIAÇÃO FND INVEST IE NE PAPEL GEO GLOBAL EQ FIA IE FND'
$ws2.Range("D423").Value = 'VARIAÇÃO FND INVEST IE NE PAPEL W V DOL MASTER FIAIE FND'
$ws2.Range("D424").Value = 'VARIAÇÃO FUNDO DE ACOES PAPEL BRA SMALL CAPS P FIA FND'
$ws2.Range("D425").Value = 'VARIAÇÃO FUNDO DE ACOES PAPEL SULAME SELECTION FIA FND'
$ws2.Range("D426").Value = 'VARIAÇÃO FUNDO DE ACOES PAPEL SULAMERICA EQUIT FIA FND'
$ws2.Range("D427").Value = 'VARIAÇÃO FUNDO INVEST BDR PAPEL SAFRA CONS AMER FIA FND'

# Leave CONTAS FAPERS parked at B9 (its last active cell) and make HP FAPERS the
# active/visible tab with D15 selected, as the workbook was left open to the user.
$ws1.Range("B9").Select()
$ws2.Activate()
$ws2.Range("D15").Select()
